# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
# Price/Volume(1h) cells are plain text (e.g. "67.888.14" w/ locale dots,
# or "  -0.96%  " with padding), so for cells whose new text parses as a
# plain number we force NumberFormat="@" (Text) before assigning, then
# reset the style back to Normal so no stray numeric formatting/quote-
# prefix styling is left behind - this keeps values stored as text just
# like the source file instead of being auto-coerced to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.888.14'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').Value = '3.267.67'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.602'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.98%  '
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.409'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.64%  '
$ws.Range('D12').Value = '3.840.09'
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.42'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.94%  '
$ws.Range('D15').Value = '67.873.39'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('E16').Value = '  -1.97%  '
$ws.Range('D17').Value = '3.296.18'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.70'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.41'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '402.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.12%  '
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.509'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.30%  '
$ws.Range('E25').Value = '  -1.68%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.66'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.47'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.86%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.89'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.32%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '164.21'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('E36').Value = '  -3.53%  '
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('E38').Value = '  +2.81%  '
$ws.Range('E39').Value = '  -3.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.48'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.34'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.07%  '
$ws.Range('D42').Value = '2.678.77'
$ws.Range('E42').Value = '  +2.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0678'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('E45').Value = '  -2.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '335.60'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.62'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('E48').Value = '  -2.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.30'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('E50').Value = '  -1.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.967'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.33%  '
